# Add data for 2021-11-04
# Update sheet title and October label to reflect new "through" date
# and update October / Total row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Through 2021-10-27"

# Update October row label (row 11, column A)
$ws.Range("A11").Value = "October (through 10-27)"

# Update October row values (row 11): B,D,E,F,G,H change; C stays the same
$ws.Range("B11").Value = 27
$ws.Range("D11").Value = 67
$ws.Range("E11").Value = 57
$ws.Range("F11").Value = 52
$ws.Range("G11").Value = 133
$ws.Range("H11").Value = 173

# Update Total row values (row 12): B,D,E,F,G,H change; C stays the same
$ws.Range("B12").Value = 253
$ws.Range("D12").Value = 694
$ws.Range("E12").Value = 605
$ws.Range("F12").Value = 474
$ws.Range("G12").Value = 1034
$ws.Range("H12").Value = 1420
